# Generate Report for Handoff
# The b.md file has been handed off again; refresh its status / handoff
# file / handoff datetime on every sheet that tracks it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is the "b.md" row (A3 = "b.md").
#   B3 (zh-cn status), C3 (de-de status) -> "Ready for handoff"
#   D3 (Latest Handoff Date)             -> "2016-08-18 07:08:19"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-08-18 07:08:19"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 is the "b.md" row.
#   C3 (Status)                  -> "Ready for handoff"
#   D3 (Latest Handoff File)     -> new handoff xlf file name
#   E3 (Latest Handoff Datetime) -> "2016-03-18 07:08:16"
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-18 07:08:16"

# Keep the hyperlink's cached display text in sync with the new file name
# (the underlying target / relationship id is unchanged). Snapshot Count
# first - mutating TextToDisplay appends a fresh entry to the collection,
# and a live upper bound would walk into it.
$zhcnLinks = $zhcn.Hyperlinks
$zhcnLinkCount = $zhcnLinks.Count
for ($i = 1; $i -le $zhcnLinkCount; $i++) {
    $link = $zhcnLinks.Item($i)
    if ($link.Range.Row -eq 3 -and $link.Range.Column -eq 4) {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet - row 3 is the "b.md" row.
#   C3 (Status)                  -> "Ready for handoff"
#   D3 (Latest Handoff File)     -> new handoff xlf file name
#   E3 (Latest Handoff Datetime) -> "2016-03-18 07:08:19"
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-18 07:08:19"

$dedeLinks = $dede.Hyperlinks
$dedeLinkCount = $dedeLinks.Count
for ($i = 1; $i -le $dedeLinkCount; $i++) {
    $link = $dedeLinks.Item($i)
    if ($link.Range.Row -eq 3 -and $link.Range.Column -eq 4) {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
